$wb = $excel.ActiveWorkbook

# AddCustomerTest sheet: normalise the "runmode" column to lowercase y/n
$ws1 = $wb.Worksheets.Item("AddCustomerTest")
$ws1.Range("E2").Value = "y"
$ws1.Range("E3").Value = "n"

# Make AddCustomerTest the active sheet/tab with E7 selected
$ws1.Range("E7").Select()
